# Updates numeric price/profit figures on several per-sheet Leve tables
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market data
# pulled in by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 999.75
$ws.Cells.Item(29, 10).Value = 999.75
$ws.Cells.Item(29, 12).Value = 2999.25
$ws.Cells.Item(29, 14).Value = -3561.25

$ws.Cells.Item(51, 8).Value = 5625.7334
$ws.Cells.Item(51, 9).Value = 6199.8887
$ws.Cells.Item(51, 11).Value = 6199.8887
$ws.Cells.Item(51, 13).Value = -5715.8887

$ws.Cells.Item(113, 8).Value = 13318.7
$ws.Cells.Item(113, 9).Value = 12915.286
$ws.Cells.Item(113, 10).Value = 14260
$ws.Cells.Item(113, 11).Value = 12915.286
$ws.Cells.Item(113, 12).Value = 14260
$ws.Cells.Item(113, 13).Value = -9661.286
$ws.Cells.Item(113, 14).Value = -20768

$ws.Cells.Item(135, 8).Value = 6497.231
$ws.Cells.Item(135, 9).Value = 7639.3
$ws.Cells.Item(135, 10).Value = 2690.3333
$ws.Cells.Item(135, 11).Value = 68753.7
$ws.Cells.Item(135, 12).Value = 24212.9997
$ws.Cells.Item(135, 13).Value = -66218.7
$ws.Cells.Item(135, 14).Value = -29282.9997

$ws.Cells.Item(138, 8).Value = 265755.03
$ws.Cells.Item(138, 10).Value = 3343.2693
$ws.Cells.Item(138, 12).Value = 10029.8079
$ws.Cells.Item(138, 14).Value = -20309.8079

$ws.Cells.Item(141, 8).Value = 11144.889
$ws.Cells.Item(141, 9).Value = 11288
$ws.Cells.Item(141, 11).Value = 33864
$ws.Cells.Item(141, 13).Value = -28684

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3707.682
$ws.Cells.Item(2, 9).Value = 3717.75
$ws.Cells.Item(2, 10).Value = 3680.8333
$ws.Cells.Item(2, 11).Value = 3717.75
$ws.Cells.Item(2, 12).Value = 3680.8333
$ws.Cells.Item(2, 13).Value = -3604.75
$ws.Cells.Item(2, 14).Value = -3906.8333

$ws.Cells.Item(32, 8).Value = 4919.722
$ws.Cells.Item(32, 9).Value = 4528.031
$ws.Cells.Item(32, 10).Value = 8556.857
$ws.Cells.Item(32, 11).Value = 4528.031
$ws.Cells.Item(32, 12).Value = 8556.857
$ws.Cells.Item(32, 13).Value = -4241.031
$ws.Cells.Item(32, 14).Value = -9130.857

$ws.Cells.Item(61, 8).Value = 9823.913
$ws.Cells.Item(61, 9).Value = 11370
$ws.Cells.Item(61, 10).Value = 4258
$ws.Cells.Item(61, 11).Value = 11370
$ws.Cells.Item(61, 12).Value = 4258
$ws.Cells.Item(61, 13).Value = -11158
$ws.Cells.Item(61, 14).Value = -4682

$ws.Cells.Item(116, 8).Value = 3707.682
$ws.Cells.Item(116, 9).Value = 3717.75
$ws.Cells.Item(116, 10).Value = 3680.8333
$ws.Cells.Item(116, 11).Value = 3717.75
$ws.Cells.Item(116, 12).Value = 3680.8333
$ws.Cells.Item(116, 13).Value = -1423.75
$ws.Cells.Item(116, 14).Value = -8268.8333

$ws.Cells.Item(132, 8).Value = 2105.6736
$ws.Cells.Item(132, 9).Value = 1485.3889
$ws.Cells.Item(132, 11).Value = 4456.1667
$ws.Cells.Item(132, 13).Value = -1926.1667

$ws.Cells.Item(136, 8).Value = 9823.913
$ws.Cells.Item(136, 9).Value = 11370
$ws.Cells.Item(136, 10).Value = 4258
$ws.Cells.Item(136, 11).Value = 34110
$ws.Cells.Item(136, 12).Value = 12774
$ws.Cells.Item(136, 13).Value = -31560
$ws.Cells.Item(136, 14).Value = -17874

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3707.682
$ws.Cells.Item(3, 9).Value = 3717.75
$ws.Cells.Item(3, 10).Value = 3680.8333
$ws.Cells.Item(3, 11).Value = 3717.75
$ws.Cells.Item(3, 12).Value = 3680.8333
$ws.Cells.Item(3, 13).Value = -3603.75
$ws.Cells.Item(3, 14).Value = -3908.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 597.3333
$ws.Cells.Item(22, 9).Value = 390.9091
$ws.Cells.Item(22, 11).Value = 390.9091
$ws.Cells.Item(22, 13).Value = -40.90910000000002

$ws.Cells.Item(31, 8).Value = 6331.759
$ws.Cells.Item(31, 9).Value = 6567.5835
$ws.Cells.Item(31, 10).Value = 5199.8
$ws.Cells.Item(31, 11).Value = 6567.5835
$ws.Cells.Item(31, 12).Value = 5199.8
$ws.Cells.Item(31, 13).Value = -6272.5835
$ws.Cells.Item(31, 14).Value = -5789.8

$ws.Cells.Item(34, 8).Value = 6331.759
$ws.Cells.Item(34, 9).Value = 6567.5835
$ws.Cells.Item(34, 10).Value = 5199.8
$ws.Cells.Item(34, 11).Value = 6567.5835
$ws.Cells.Item(34, 12).Value = 5199.8
$ws.Cells.Item(34, 13).Value = -6365.5835
$ws.Cells.Item(34, 14).Value = -5603.8

$ws.Cells.Item(99, 8).Value = 336353.2
$ws.Cells.Item(99, 9).Value = 627281.1
$ws.Cells.Item(99, 10).Value = 3864.1428
$ws.Cells.Item(99, 11).Value = 627281.1
$ws.Cells.Item(99, 12).Value = 3864.1428
$ws.Cells.Item(99, 13).Value = -625783.1
$ws.Cells.Item(99, 14).Value = -6860.1428

$ws.Cells.Item(126, 8).Value = 336353.2
$ws.Cells.Item(126, 9).Value = 627281.1
$ws.Cells.Item(126, 10).Value = 3864.1428
$ws.Cells.Item(126, 11).Value = 1881843.3
$ws.Cells.Item(126, 12).Value = 11592.4284
$ws.Cells.Item(126, 13).Value = -1879373.3
$ws.Cells.Item(126, 14).Value = -16532.4284

$ws.Cells.Item(134, 8).Value = 7741.25
$ws.Cells.Item(134, 9).Value = 8057.4
$ws.Cells.Item(134, 11).Value = 24172.2
$ws.Cells.Item(134, 13).Value = -21637.2

$ws.Cells.Item(141, 8).Value = 291471.34
$ws.Cells.Item(141, 10).Value = 320405.25
$ws.Cells.Item(141, 12).Value = 320405.25
$ws.Cells.Item(141, 14).Value = -330765.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 37.42857
$ws.Cells.Item(12, 10).Value = 41.57143
$ws.Cells.Item(12, 12).Value = 124.71429
$ws.Cells.Item(12, 14).Value = -470.71429

$ws.Cells.Item(68, 8).Value = 10302.267
$ws.Cells.Item(68, 10).Value = 12423.667
$ws.Cells.Item(68, 12).Value = 37271.001
$ws.Cells.Item(68, 14).Value = -38893.001

$ws.Cells.Item(71, 8).Value = 10302.267
$ws.Cells.Item(71, 10).Value = 12423.667
$ws.Cells.Item(71, 12).Value = 111813.003
$ws.Cells.Item(71, 14).Value = -119925.003

$ws.Cells.Item(80, 8).Value = 64689.285
$ws.Cells.Item(80, 9).Value = 6999
$ws.Cells.Item(80, 10).Value = 69127
$ws.Cells.Item(80, 11).Value = 20997
$ws.Cells.Item(80, 12).Value = 207381
$ws.Cells.Item(80, 13).Value = -20061
$ws.Cells.Item(80, 14).Value = -209253

$ws.Cells.Item(83, 8).Value = 64689.285
$ws.Cells.Item(83, 9).Value = 6999
$ws.Cells.Item(83, 10).Value = 69127
$ws.Cells.Item(83, 11).Value = 62991
$ws.Cells.Item(83, 12).Value = 622143
$ws.Cells.Item(83, 13).Value = -58311
$ws.Cells.Item(83, 14).Value = -631503

$ws.Cells.Item(122, 8).Value = 5037
$ws.Cells.Item(122, 9).Value = 1475.7142
$ws.Cells.Item(122, 10).Value = 5896.6206
$ws.Cells.Item(122, 11).Value = 13281.4278
$ws.Cells.Item(122, 12).Value = 53069.5854
$ws.Cells.Item(122, 13).Value = -10831.4278
$ws.Cells.Item(122, 14).Value = -57969.5854

$ws.Cells.Item(131, 8).Value = 1936.4382
$ws.Cells.Item(131, 9).Value = 1199.7142
$ws.Cells.Item(131, 11).Value = 3599.1426
$ws.Cells.Item(131, 13).Value = 1440.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 15950.9
$ws.Cells.Item(102, 9).Value = 17056.555
$ws.Cells.Item(102, 11).Value = 17056.555
$ws.Cells.Item(102, 13).Value = -15434.555

$ws.Cells.Item(122, 8).Value = 6162.081
$ws.Cells.Item(122, 9).Value = 4057.0571
$ws.Cells.Item(122, 10).Value = 43000
$ws.Cells.Item(122, 11).Value = 12171.1713
$ws.Cells.Item(122, 12).Value = 129000
$ws.Cells.Item(122, 13).Value = -9721.1713
$ws.Cells.Item(122, 14).Value = -133900

$ws.Cells.Item(126, 8).Value = 5414.452
$ws.Cells.Item(126, 9).Value = 6226.241
$ws.Cells.Item(126, 10).Value = 3603.5386
$ws.Cells.Item(126, 11).Value = 18678.723
$ws.Cells.Item(126, 12).Value = 10810.6158
$ws.Cells.Item(126, 13).Value = -16208.723
$ws.Cells.Item(126, 14).Value = -15750.6158

$ws.Cells.Item(133, 8).Value = 69995
$ws.Cells.Item(133, 10).Value = 69995
$ws.Cells.Item(133, 12).Value = 69995
$ws.Cells.Item(133, 14).Value = -80115

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 22210.545
$ws.Cells.Item(7, 9).Value = 28883.375
$ws.Cells.Item(7, 10).Value = 4416.3335
$ws.Cells.Item(7, 11).Value = 28883.375
$ws.Cells.Item(7, 12).Value = 4416.3335
$ws.Cells.Item(7, 13).Value = -28771.375
$ws.Cells.Item(7, 14).Value = -4640.3335

$ws.Cells.Item(22, 8).Value = 10172.728
$ws.Cells.Item(22, 9).Value = 14885.714
$ws.Cells.Item(22, 11).Value = 14885.714
$ws.Cells.Item(22, 13).Value = -14590.714

$ws.Cells.Item(27, 8).Value = 10172.728
$ws.Cells.Item(27, 9).Value = 14885.714
$ws.Cells.Item(27, 11).Value = 14885.714
$ws.Cells.Item(27, 13).Value = -14778.714

$ws.Cells.Item(35, 8).Value = 832.1667
$ws.Cells.Item(35, 9).Value = 878.6
$ws.Cells.Item(35, 10).Value = 600
$ws.Cells.Item(35, 11).Value = 878.6
$ws.Cells.Item(35, 12).Value = 600
$ws.Cells.Item(35, 13).Value = -542.6
$ws.Cells.Item(35, 14).Value = -1272

$ws.Cells.Item(122, 8).Value = 7032.846
$ws.Cells.Item(122, 9).Value = 7643.3
$ws.Cells.Item(122, 11).Value = 22929.9
$ws.Cells.Item(122, 13).Value = -20479.9

$ws.Cells.Item(126, 8).Value = 22210.545
$ws.Cells.Item(126, 9).Value = 28883.375
$ws.Cells.Item(126, 10).Value = 4416.3335
$ws.Cells.Item(126, 11).Value = 86650.125
$ws.Cells.Item(126, 12).Value = 13249.0005
$ws.Cells.Item(126, 13).Value = -84180.125
$ws.Cells.Item(126, 14).Value = -18189.0005

$ws.Cells.Item(136, 8).Value = 4676.15
$ws.Cells.Item(136, 9).Value = 2780.2856
$ws.Cells.Item(136, 11).Value = 8340.856800000001
$ws.Cells.Item(136, 13).Value = -5790.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 20000
$ws.Cells.Item(40, 10).Value = 20000
$ws.Cells.Item(40, 12).Value = 20000
$ws.Cells.Item(40, 14).Value = -20298

$ws.Cells.Item(62, 8).Value = 347919.9
$ws.Cells.Item(62, 9).Value = 491428.44
$ws.Cells.Item(62, 10).Value = 13066.667
$ws.Cells.Item(62, 11).Value = 491428.44
$ws.Cells.Item(62, 12).Value = 13066.667
$ws.Cells.Item(62, 13).Value = -490804.44
$ws.Cells.Item(62, 14).Value = -14314.667

$ws.Cells.Item(65, 8).Value = 347919.9
$ws.Cells.Item(65, 9).Value = 491428.44
$ws.Cells.Item(65, 10).Value = 13066.667
$ws.Cells.Item(65, 11).Value = 2457142.2
$ws.Cells.Item(65, 12).Value = 65333.335
$ws.Cells.Item(65, 13).Value = -2454022.2
$ws.Cells.Item(65, 14).Value = -71573.33499999999

$ws.Cells.Item(126, 8).Value = 17611.793
$ws.Cells.Item(126, 9).Value = 20456.25
$ws.Cells.Item(126, 10).Value = 3958.4
$ws.Cells.Item(126, 11).Value = 61368.75
$ws.Cells.Item(126, 12).Value = 11875.2
$ws.Cells.Item(126, 13).Value = -58898.75
$ws.Cells.Item(126, 14).Value = -16815.2

$ws.Cells.Item(136, 8).Value = 434298.94
$ws.Cells.Item(136, 9).Value = 456816.53
$ws.Cells.Item(136, 11).Value = 1370449.59
$ws.Cells.Item(136, 13).Value = -1367899.59

$ws.Cells.Item(139, 8).Value = 115666.336
$ws.Cells.Item(139, 10).Value = 68499.5
$ws.Cells.Item(139, 12).Value = 68499.5
$ws.Cells.Item(139, 14).Value = -78779.5

$ws.Cells.Item(141, 8).Value = 94665.836
$ws.Cells.Item(141, 10).Value = 83799
$ws.Cells.Item(141, 12).Value = 83799
$ws.Cells.Item(141, 14).Value = -94159
